$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 4): Divisor / Block Size / Ratio
$ws.Range("A4").Value = "Divisor"
$ws.Range("B4").Value = "Block Size"
$ws.Range("C4").Value = "Ratio "

# Column B width to fit the "Block Size" header (bestFit column, final width 12)
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666

# Update selection to match the diff (F7:F8, active cell F7)
$ws.Range("F7:F8").Select()
